# Add a new time-log entry on Sheet1 row 99 (date 2014-10-24, 00:00 - 01:00, no
# interruption, Activity = Coding), mirroring the existing pattern of rows in the
# log. Sheet2 / the chart derive their values from Sheet1 via formulas and will
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New entry: Date, Start Time, Stop Time, Interruption (mins), Delta (formula
# already present), Activity.
$ws1.Cells.Item(99, 1).Value = 41936
$ws1.Cells.Item(99, 2).Value = 0
$ws1.Cells.Item(99, 3).Value = (1.0 / 24.0)
$ws1.Cells.Item(99, 4).Value = 0
$ws1.Cells.Item(99, 6).Value = "Coding"

# Update the active selection to A100, as it would be after entering the row.
$ws1.Range("A100").Select()

$excel.Calculate()
